# Trade #18 closed at 2026-02-17 20:53:04 - unknown UNKNOWN +0.000%
#
# Updates the live trading results workbook:
#  - Summary: roll up capital / P&L / trade counters
#  - Strategy Status: roll up the MarketMaking strategy row
#  - All Trades: close out the existing OPEN MarketMaking trade (row for
#    Trade #46) and append the newly opened Trade #79
#  - MarketMaking: same two updates, mirrored onto the per-strategy sheet
#    (column order differs from "All Trades" on this sheet)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1400.41   # Current Capital
$summary.Range("B4").Value = 0.2       # Total P&L $
$summary.Range("B5").Value = 0.09      # Total P&L %
$summary.Range("B6").Value = 46        # Total Trades
$summary.Range("B7").Value = 21        # Winning Trades
$summary.Range("B9").Value = 45.65     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status (MarketMaking row = row 5)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 100.41
$status.Range("D5").Value = 13
$status.Range("E5").Value = 0.09
$status.Range("F5").Value = 0.41
$status.Range("G5").Value = 53.85

# ---------------------------------------------------------------------
# All Trades
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Close out Trade #46 (row 47): was OPEN, now CLOSED via early_exit
$allTrades.Range("G47").Value = 0.135922
$allTrades.Range("H47").Value = "CLOSED"
$allTrades.Range("I47").Value = 4.5556
$allTrades.Range("J47").Value = 0.01
$allTrades.Range("K47").Value = 100.41
$allTrades.Range("L47").Value = "early_exit"
$allTrades.Range("M47").Value = 0.14

# Append newly opened Trade #79 as row 80
$allTrades.Range("A80").Value = 79
$allTrades.Range("B80").NumberFormat = "@"
$allTrades.Range("B80").Value = "2026-02-17"
$allTrades.Range("C80").NumberFormat = "@"
$allTrades.Range("C80").Value = "20:52:58"
$allTrades.Range("D80").Value = "MarketMaking"
$allTrades.Range("E80").Value = "UP"
$allTrades.Range("F80").Value = 0.13
$allTrades.Range("H80").Value = "OPEN"
$allTrades.Range("I80").Value = 0
$allTrades.Range("J80").Value = 0
$allTrades.Range("K80").Value = 100.4010404544114
$allTrades.Range("M80").Value = 0
$allTrades.Range("N80").Value = 0
$allTrades.Range("O80").Value = 0
$allTrades.Range("P80").Value = 0.6
$allTrades.Range("Q80").Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------
# MarketMaking (per-strategy sheet; note L/M/N/O/P/Q order differs from
# "All Trades": L/M = slippage, N = confidence, O = entry reason,
# P = exit reason, Q = duration)
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

# Close out Trade #46 (row 14)
$mm.Range("G14").Value = 0.135922
$mm.Range("H14").Value = "CLOSED"
$mm.Range("I14").Value = 4.5556
$mm.Range("J14").Value = 0.01
$mm.Range("K14").Value = 100.41
$mm.Range("P14").Value = "early_exit"
$mm.Range("Q14").Value = 0.14

# Append newly opened Trade #79 as row 47
$mm.Range("A47").Value = 79
$mm.Range("B47").NumberFormat = "@"
$mm.Range("B47").Value = "2026-02-17"
$mm.Range("C47").NumberFormat = "@"
$mm.Range("C47").Value = "20:52:58"
$mm.Range("D47").Value = "MarketMaking"
$mm.Range("E47").Value = "UP"
$mm.Range("F47").Value = 0.13
$mm.Range("H47").Value = "OPEN"
$mm.Range("I47").Value = 0
$mm.Range("J47").Value = 0
$mm.Range("K47").Value = 100.4010404544114
$mm.Range("L47").Value = 0
$mm.Range("M47").Value = 0
$mm.Range("N47").Value = 0.6
$mm.Range("O47").Value = "Normal spread capture: 19600 bps"
$mm.Range("Q47").Value = 0

Write-Output "edit.ps1 applied"
